$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New participant rows (sub_031 .. sub_044), appended after existing row 31.
$data = @(
    @{ Row = 32; Participant = "sub_031"; Fail = $false; Reason = $null },
    @{ Row = 33; Participant = "sub_032"; Fail = $true;  Reason = "scroll bar issues" },
    @{ Row = 34; Participant = "sub_033"; Fail = $true;  Reason = "didn’t list all the visible items" },
    @{ Row = 35; Participant = "sub_034"; Fail = $false; Reason = $null },
    @{ Row = 36; Participant = "sub_035"; Fail = $false; Reason = $null },
    @{ Row = 37; Participant = "sub_036"; Fail = $false; Reason = $null },
    @{ Row = 38; Participant = "sub_037"; Fail = $false; Reason = $null },
    @{ Row = 39; Participant = "sub_038"; Fail = $true;  Reason = "didn’t list all the items" },
    @{ Row = 40; Participant = "sub_039"; Fail = $false; Reason = $null },
    @{ Row = 41; Participant = "sub_040"; Fail = $false; Reason = $null },
    @{ Row = 42; Participant = "sub_041"; Fail = $false; Reason = $null },
    @{ Row = 43; Participant = "sub_042"; Fail = $false; Reason = $null },
    @{ Row = 44; Participant = "sub_043"; Fail = $false; Reason = $null },
    @{ Row = 45; Participant = "sub_044"; Fail = $false; Reason = $null }
)

# Populate column A (participant) for all new rows first so the new
# shared-string entries are interned in participant order, matching how
# the source workbook was produced (sub_031..sub_044 before the reasons).
foreach ($entry in $data) {
    $ws.Cells.Item($entry.Row, 1).Value = $entry.Participant
}

foreach ($entry in $data) {
    $ws.Cells.Item($entry.Row, 2).Value = $entry.Fail
}

foreach ($entry in $data) {
    if ($entry.Reason) {
        $ws.Cells.Item($entry.Row, 3).Value = $entry.Reason
    }
}

# Update view to reflect newly scrolled position / selection as in the saved file.
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Range("C40").Select() | Out-Null
